$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear out the old string-backed columns (A-D) for rows 2-7 so the
# shared-string table drops refs to the old entries before we rewrite
# everything (rows 2-10) with the updated cluster/gene labels.
$ws.Range("A2:D7").ClearContents()

# Row 2
$ws.Cells.Item(2,1).Value = "ECs"
$ws.Cells.Item(2,2).Value = "Tgfb2"
$ws.Cells.Item(2,3).Value = "Acvr1"
$ws.Cells.Item(2,4).Value = "ECs"
$ws.Cells.Item(2,5).Value = 2
$ws.Cells.Item(2,6).Value = 0.6666666666666666
$ws.Cells.Item(2,7).Value = 1.465615333333333
$ws.Cells.Item(2,8).Value = 4.396846
$ws.Cells.Item(2,9).Value = 0.04672291954663727
$ws.Cells.Item(2,10).Value = 0.04672291954663728
$ws.Cells.Item(2,11).Value = 3
$ws.Cells.Item(2,12).Value = 1
$ws.Cells.Item(2,13).Value = 4.695610666666666
$ws.Cells.Item(2,14).Value = 14.086832
$ws.Cells.Item(2,15).Value = 0.1802066564018305
$ws.Cells.Item(2,16).Value = 0.1802066564018305
$ws.Cells.Item(2,17).Value = 6.881958992430222
$ws.Cells.Item(2,18).Value = 61.937630931872
$ws.Cells.Item(2,19).Value = 0.008419781108831233
$ws.Cells.Item(2,20).Value = 0.008419781108831233

# Row 3
$ws.Cells.Item(3,1).Value = "ECs"
$ws.Cells.Item(3,2).Value = "Tgfb2"
$ws.Cells.Item(3,3).Value = "Acvr1"
$ws.Cells.Item(3,4).Value = "FAPs"
$ws.Cells.Item(3,5).Value = 2
$ws.Cells.Item(3,6).Value = 0.6666666666666666
$ws.Cells.Item(3,7).Value = 1.465615333333333
$ws.Cells.Item(3,8).Value = 4.396846
$ws.Cells.Item(3,9).Value = 0.04672291954663727
$ws.Cells.Item(3,10).Value = 0.04672291954663728
$ws.Cells.Item(3,11).Value = 3
$ws.Cells.Item(3,12).Value = 1
$ws.Cells.Item(3,13).Value = 15.51448033333333
$ws.Cells.Item(3,14).Value = 46.543441
$ws.Cells.Item(3,15).Value = 0.5954098039960916
$ws.Cells.Item(3,16).Value = 0.5954098039960916
$ws.Cells.Item(3,17).Value = 22.73826026523178
$ws.Cells.Item(3,18).Value = 204.644342387086
$ws.Cells.Item(3,19).Value = 0.02781928436938845
$ws.Cells.Item(3,20).Value = 0.02781928436938846

# Row 4
$ws.Cells.Item(4,1).Value = "ECs"
$ws.Cells.Item(4,2).Value = "Tgfb2"
$ws.Cells.Item(4,3).Value = "Acvr1"
$ws.Cells.Item(4,4).Value = "sCs"
$ws.Cells.Item(4,5).Value = 2
$ws.Cells.Item(4,6).Value = 0.6666666666666666
$ws.Cells.Item(4,7).Value = 1.465615333333333
$ws.Cells.Item(4,8).Value = 4.396846
$ws.Cells.Item(4,9).Value = 0.04672291954663727
$ws.Cells.Item(4,10).Value = 0.04672291954663728
$ws.Cells.Item(4,11).Value = 3
$ws.Cells.Item(4,12).Value = 1
$ws.Cells.Item(4,13).Value = 5.846719333333333
$ws.Cells.Item(4,14).Value = 17.540158
$ws.Cells.Item(4,15).Value = 0.2243835396020779
$ws.Cells.Item(4,16).Value = 0.2243835396020779
$ws.Cells.Item(4,17).Value = 8.569041504629777
$ws.Cells.Item(4,18).Value = 77.12137354166799
$ws.Cells.Item(4,19).Value = 0.01048385406841758
$ws.Cells.Item(4,20).Value = 0.01048385406841758

# Row 5
$ws.Cells.Item(5,1).Value = "FAPs"
$ws.Cells.Item(5,2).Value = "Tgfb2"
$ws.Cells.Item(5,3).Value = "Acvr1"
$ws.Cells.Item(5,4).Value = "ECs"
$ws.Cells.Item(5,5).Value = 3
$ws.Cells.Item(5,6).Value = 1
$ws.Cells.Item(5,7).Value = 18.88237266666667
$ws.Cells.Item(5,8).Value = 56.64711800000001
$ws.Cells.Item(5,9).Value = 0.6019584804341267
$ws.Cells.Item(5,10).Value = 0.6019584804341268
$ws.Cells.Item(5,11).Value = 3
$ws.Cells.Item(5,12).Value = 1
$ws.Cells.Item(5,13).Value = 4.695610666666666
$ws.Cells.Item(5,14).Value = 14.086832
$ws.Cells.Item(5,15).Value = 0.1802066564018305
$ws.Cells.Item(5,16).Value = 0.1802066564018305
$ws.Cells.Item(5,17).Value = 88.66427050557512
$ws.Cells.Item(5,18).Value = 797.9784345501761
$ws.Cells.Item(5,19).Value = 0.1084769250517607
$ws.Cells.Item(5,20).Value = 0.1084769250517607

# Row 6
$ws.Cells.Item(6,1).Value = "FAPs"
$ws.Cells.Item(6,2).Value = "Tgfb2"
$ws.Cells.Item(6,3).Value = "Acvr1"
$ws.Cells.Item(6,4).Value = "FAPs"
$ws.Cells.Item(6,5).Value = 3
$ws.Cells.Item(6,6).Value = 1
$ws.Cells.Item(6,7).Value = 18.88237266666667
$ws.Cells.Item(6,8).Value = 56.64711800000001
$ws.Cells.Item(6,9).Value = 0.6019584804341267
$ws.Cells.Item(6,10).Value = 0.6019584804341268
$ws.Cells.Item(6,11).Value = 3
$ws.Cells.Item(6,12).Value = 1
$ws.Cells.Item(6,13).Value = 15.51448033333333
$ws.Cells.Item(6,14).Value = 46.543441
$ws.Cells.Item(6,15).Value = 0.5954098039960916
$ws.Cells.Item(6,16).Value = 0.5954098039960916
$ws.Cells.Item(6,17).Value = 292.9501993836709
$ws.Cells.Item(6,18).Value = 2636.551794453038
$ws.Cells.Item(6,19).Value = 0.3584119808490685
$ws.Cells.Item(6,20).Value = 0.3584119808490686

# Row 7
$ws.Cells.Item(7,1).Value = "FAPs"
$ws.Cells.Item(7,2).Value = "Tgfb2"
$ws.Cells.Item(7,3).Value = "Acvr1"
$ws.Cells.Item(7,4).Value = "sCs"
$ws.Cells.Item(7,5).Value = 3
$ws.Cells.Item(7,6).Value = 1
$ws.Cells.Item(7,7).Value = 18.88237266666667
$ws.Cells.Item(7,8).Value = 56.64711800000001
$ws.Cells.Item(7,9).Value = 0.6019584804341267
$ws.Cells.Item(7,10).Value = 0.6019584804341268
$ws.Cells.Item(7,11).Value = 3
$ws.Cells.Item(7,12).Value = 1
$ws.Cells.Item(7,13).Value = 5.846719333333333
$ws.Cells.Item(7,14).Value = 17.540158
$ws.Cells.Item(7,15).Value = 0.2243835396020779
$ws.Cells.Item(7,16).Value = 0.2243835396020779
$ws.Cells.Item(7,17).Value = 110.3999333294049
$ws.Cells.Item(7,18).Value = 993.599399964644
$ws.Cells.Item(7,19).Value = 0.1350695745332975
$ws.Cells.Item(7,20).Value = 0.1350695745332975

# Row 8
$ws.Cells.Item(8,1).Value = "sCs"
$ws.Cells.Item(8,2).Value = "Tgfb2"
$ws.Cells.Item(8,3).Value = "Acvr1"
$ws.Cells.Item(8,4).Value = "ECs"
$ws.Cells.Item(8,5).Value = 3
$ws.Cells.Item(8,6).Value = 1
$ws.Cells.Item(8,7).Value = 11.020243
$ws.Cells.Item(8,8).Value = 33.060729
$ws.Cells.Item(8,9).Value = 0.351318600019236
$ws.Cells.Item(8,10).Value = 0.351318600019236
$ws.Cells.Item(8,11).Value = 3
$ws.Cells.Item(8,12).Value = 1
$ws.Cells.Item(8,13).Value = 4.695610666666666
$ws.Cells.Item(8,14).Value = 14.086832
$ws.Cells.Item(8,15).Value = 0.1802066564018305
$ws.Cells.Item(8,16).Value = 0.1802066564018305
$ws.Cells.Item(8,17).Value = 51.74677058005867
$ws.Cells.Item(8,18).Value = 465.720935220528
$ws.Cells.Item(8,19).Value = 0.06330995024123859
$ws.Cells.Item(8,20).Value = 0.06330995024123859

# Row 9
$ws.Cells.Item(9,1).Value = "sCs"
$ws.Cells.Item(9,2).Value = "Tgfb2"
$ws.Cells.Item(9,3).Value = "Acvr1"
$ws.Cells.Item(9,4).Value = "FAPs"
$ws.Cells.Item(9,5).Value = 3
$ws.Cells.Item(9,6).Value = 1
$ws.Cells.Item(9,7).Value = 11.020243
$ws.Cells.Item(9,8).Value = 33.060729
$ws.Cells.Item(9,9).Value = 0.351318600019236
$ws.Cells.Item(9,10).Value = 0.351318600019236
$ws.Cells.Item(9,11).Value = 3
$ws.Cells.Item(9,12).Value = 1
$ws.Cells.Item(9,13).Value = 15.51448033333333
$ws.Cells.Item(9,14).Value = 46.543441
$ws.Cells.Item(9,15).Value = 0.5954098039960916
$ws.Cells.Item(9,16).Value = 0.5954098039960916
$ws.Cells.Item(9,17).Value = 170.9733432920544
$ws.Cells.Item(9,18).Value = 1538.760089628489
$ws.Cells.Item(9,19).Value = 0.2091785387776346
$ws.Cells.Item(9,20).Value = 0.2091785387776346

# Row 10
$ws.Cells.Item(10,1).Value = "sCs"
$ws.Cells.Item(10,2).Value = "Tgfb2"
$ws.Cells.Item(10,3).Value = "Acvr1"
$ws.Cells.Item(10,4).Value = "sCs"
$ws.Cells.Item(10,5).Value = 3
$ws.Cells.Item(10,6).Value = 1
$ws.Cells.Item(10,7).Value = 11.020243
$ws.Cells.Item(10,8).Value = 33.060729
$ws.Cells.Item(10,9).Value = 0.351318600019236
$ws.Cells.Item(10,10).Value = 0.351318600019236
$ws.Cells.Item(10,11).Value = 3
$ws.Cells.Item(10,12).Value = 1
$ws.Cells.Item(10,13).Value = 5.846719333333333
$ws.Cells.Item(10,14).Value = 17.540158
$ws.Cells.Item(10,15).Value = 0.2243835396020779
$ws.Cells.Item(10,16).Value = 0.2243835396020779
$ws.Cells.Item(10,17).Value = 64.43226780613134
$ws.Cells.Item(10,18).Value = 579.890410255182
$ws.Cells.Item(10,19).Value = 0.07883011100036282
$ws.Cells.Item(10,20).Value = 0.07883011100036282
